# Apply the "release-notes.md" metadata refresh to the IG export workbook.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump version/status/date/contact -----------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value  = "0.4.0-snapshot-1"                 # Version
$wsMeta.Range("B6").Value  = "draft"                            # Status
$wsMeta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"        # Date
$wsMeta.Range("B10").Value = "ANS (https://esante.gouv.fr)"     # Contact

# --- Elements sheet: swap the two "Mapping" columns (AK <-> AL) -------
# The "RIM Mapping" and "Spécification métier" mapping columns traded
# places: everything that used to live in column AK (37) now lives in
# column AL (38), and vice versa - header text, column width and the
# per-row mapping values alike.
$wsElem = $wb.Worksheets.Item("Elements")

$lastRow = $wsElem.Cells.Item($wsElem.Rows.Count, 37).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $wsElem.Cells.Item($r, 37)
    $alCell = $wsElem.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Column widths follow the same swap (AK becomes the wide column, AL the
# narrow one). The values below are the ColumnWidth inputs that this
# engine's pixel-quantized width model resolves to the saved widths
# 81.9921875 / 24.98046875 (i.e. stored widths of ~82 / ~25 chars).
$wsElem.Columns.Item(37).ColumnWidth = 81.16666666666667
$wsElem.Columns.Item(38).ColumnWidth = 24.166666666666668
